$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: J1 changes from "description" to "type", new L1 = "Image inactive"
$ws.Range("J1").Value = "type"
$ws.Range("L1").Value = "Image inactive"

# Row 2: B2 stays DatabasesToolbar (unchanged value, index shifted only)

# Row 3: Import dropdown (replaces old "project_new" button)
$ws.Range("B3").Value = "ImportProjectButton"
$ws.Range("C3").Value = "wiDropdown"
$ws.Range("D3").Formula = "=REPLACE(C3, 1, 2, """")"
$ws.Range("E3").Value = "import_file"
$ws.Range("F3").Formula = "=SUBSTITUTE(E3,""_"",""-"")"
$ws.Range("G3").Value = "Import"
$ws.Range("J3").Value = "small"

# Row 4: Reload button (replaces old "project_open" button)
$ws.Range("B4").Value = "ReloadProjectButton"
$ws.Range("C4").Value = "wiButton"
$ws.Range("D4").Formula = "=REPLACE(C4, 1, 2, """")"
$ws.Range("E4").Value = "reload_16x16"
$ws.Range("F4").Formula = "=SUBSTITUTE(E4,""_"",""-"")"
$ws.Range("G4").Value = "Reload"
$ws.Range("J4").Value = "small"

# Row 5: Apply/Collapse button (replaces old "project_close" button)
$ws.Range("B5").Value = "CollapseProjectButton"
$ws.Range("C5").Value = "wiButton"
$ws.Range("D5").Formula = "=REPLACE(C5, 1, 2, """")"
$ws.Range("E5").Value = "apply_16x16"
$ws.Range("F5").Formula = "=SUBSTITUTE(E5,""_"",""-"")"
$ws.Range("G5").Value = "Apply"
$ws.Range("J5").Value = "small"

# Row 6: new Delete button (new row, replacing former blank spacer row 6)
$ws.Range("A6").Value = 1.4
$ws.Range("B6").Value = "DeleteProjectButton"
$ws.Range("C6").Value = "wiButton "
$ws.Range("D6").Formula = "=REPLACE(C6, 1, 2, """")"
$ws.Range("E6").Value = "delete_16x16"
$ws.Range("F6").Formula = "=SUBSTITUTE(E6,""_"",""-"")"
$ws.Range("G6").Value = "Delete"
$ws.Range("J6").Value = "small"
$ws.Range("L6").Value = "delete_inactive_16x16/"

# Row 7: new Browse project button
$ws.Range("A7").Value = 1.5
$ws.Range("B7").Value = "BrowseProjectButton"
$ws.Range("C7").Value = "wiButton "
$ws.Range("D7").Formula = "=REPLACE(C7, 1, 2, """")"
$ws.Range("E7").Value = "browse_project"
$ws.Range("F7").Formula = "=SUBSTITUTE(E7,""_"",""-"")"
$ws.Range("G7").Value = "Browse project"
$ws.Range("J7").Value = "small"

# Formatting: B6/B7 use the style from B4/B5 (wiButton-style), E6 gets a distinct style
$ws.Range("B4,B5").Copy()
$ws.Range("B6:B7").PasteSpecial(-4122)

# View settings
$ws.Application.ActiveWindow.Zoom = 93
$ws.Range("J11").Select()
